# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data snapshot from 18 Jun 2020 19:12 -> 20:29.
# This updates case counters for a number of countries and, because the
# sheet is kept sorted by "Casos totales" (column B) descending, a handful
# of countries swap rank/row with their neighbours as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 20:29"

# --- Simple data refreshes (country/row unchanged) ---

# Estados Unidos
$ws.Cells.Item(4, 2).Value = 2247521
$ws.Cells.Item(4, 3).Value = 13050
$ws.Cells.Item(4, 4).Value = 920921
$ws.Cells.Item(4, 5).Value = 1206344
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 315
$ws.Cells.Item(4, 8).Value = 120256

# India
$ws.Cells.Item(7, 2).Value = 378171
$ws.Cells.Item(7, 3).Value = 10907
$ws.Cells.Item(7, 4).Value = 201297
$ws.Cells.Item(7, 5).Value = 164335
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 277
$ws.Cells.Item(7, 8).Value = 12539

# Chile
$ws.Cells.Item(12, 2).Value = 225103
$ws.Cells.Item(12, 3).Value = 4475
$ws.Cells.Item(12, 4).Value = 186441
$ws.Cells.Item(12, 5).Value = 34821
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 226
$ws.Cells.Item(12, 8).Value = 3841

# Turquia
$ws.Cells.Item(15, 2).Value = 184031
$ws.Cells.Item(15, 3).Value = 1304
$ws.Cells.Item(15, 4).Value = 156022
$ws.Cells.Item(15, 5).Value = 23127
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 21
$ws.Cells.Item(15, 8).Value = 4882

# Francia
$ws.Cells.Item(18, 2).Value = 158641
$ws.Cells.Item(18, 3).Value = 467
$ws.Cells.Item(18, 4).Value = 73887
$ws.Cells.Item(18, 5).Value = 55151
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 28
$ws.Cells.Item(18, 8).Value = 29603

# Irlanda
$ws.Cells.Item(45, 2).Value = 25355
$ws.Cells.Item(45, 3).Value = 14
$ws.Cells.Item(45, 4).Value = 22698
$ws.Cells.Item(45, 5).Value = 943
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 1714

# Marruecos
$ws.Cells.Item(68, 2).Value = 9074
$ws.Cells.Item(68, 3).Value = 77
$ws.Cells.Item(68, 4).Value = 8041
$ws.Cells.Item(68, 5).Value = 820
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 213

# Republica de Africa Central
$ws.Cells.Item(97, 2).Value = 2605
$ws.Cells.Item(97, 3).Value = 41
$ws.Cells.Item(97, 4).Value = 417
$ws.Cells.Item(97, 5).Value = 2169
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 1
$ws.Cells.Item(97, 8).Value = 19

# Sri Lanka
$ws.Cells.Item(104, 2).Value = 1928
$ws.Cells.Item(104, 3).Value = 4
$ws.Cells.Item(104, 4).Value = 1421
$ws.Cells.Item(104, 5).Value = 496
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 11

# Yemen
$ws.Cells.Item(129, 2).Value = 909
$ws.Cells.Item(129, 3).Value = 7
$ws.Cells.Item(129, 4).Value = 273
$ws.Cells.Item(129, 5).Value = 388
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 4
$ws.Cells.Item(129, 8).Value = 248

# Malaui
$ws.Cells.Item(147, 2).Value = 592
$ws.Cells.Item(147, 3).Value = 20
$ws.Cells.Item(147, 4).Value = 74
$ws.Cells.Item(147, 5).Value = 510
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 2
$ws.Cells.Item(147, 8).Value = 8

# Suazilandia
$ws.Cells.Item(148, 2).Value = 586
$ws.Cells.Item(148, 3).Value = 23
$ws.Cells.Item(148, 4).Value = 267
$ws.Cells.Item(148, 5).Value = 315
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 4

# Surinam
$ws.Cells.Item(161, 2).Value = 261
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 74
$ws.Cells.Item(161, 5).Value = 180
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 7

# --- Rank swaps (country name changes because rows re-sort by Casos totales) ---

# rows 107-109: Sudan del Sur moves above Nicaragua / Islandia
$ws.Cells.Item(107, 1).Value = "Sudan del Sur"
$ws.Cells.Item(107, 2).Value = 1830
$ws.Cells.Item(107, 3).Value = 17
$ws.Cells.Item(107, 4).Value = 117
$ws.Cells.Item(107, 5).Value = 1681
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 32

$ws.Cells.Item(108, 1).Value = "Nicaragua"
$ws.Cells.Item(108, 2).Value = 1823
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 1238
$ws.Cells.Item(108, 5).Value = 521
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 64

$ws.Cells.Item(109, 1).Value = "Islandia"
$ws.Cells.Item(109, 2).Value = 1816
$ws.Cells.Item(109, 3).Value = 1
$ws.Cells.Item(109, 4).Value = 1801
$ws.Cells.Item(109, 5).Value = 5
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 10

# rows 168-169: Guyana moves above Gibraltar
$ws.Cells.Item(168, 1).Value = "Guyana"
$ws.Cells.Item(168, 2).Value = 183
$ws.Cells.Item(168, 3).Value = 12
$ws.Cells.Item(168, 4).Value = 102
$ws.Cells.Item(168, 5).Value = 69
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 12

$ws.Cells.Item(169, 1).Value = "Gibraltar"
$ws.Cells.Item(169, 2).Value = 176
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 176
$ws.Cells.Item(169, 5).Value = 0
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

# rows 173-174: Eritrea moves above Brunei
$ws.Cells.Item(173, 1).Value = "Eritrea"
$ws.Cells.Item(173, 2).Value = 142
$ws.Cells.Item(173, 3).Value = 11
$ws.Cells.Item(173, 4).Value = 39
$ws.Cells.Item(173, 5).Value = 103
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = "Brunei"
$ws.Cells.Item(174, 2).Value = 141
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 138
$ws.Cells.Item(174, 5).Value = 0
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 3

# rows 208-209: Santa Sede moves above Islas Turcas y Caicos
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 1

# rows 213-214: Papua Nueva Guinea moves above Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1
